# Advanced OCR v2.0: Image preprocessing + item extraction
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet
$ws.Name = "Delivery Notes"

# Insert a new column F ("Items Count") before the existing "Review Status" column.
$ws.Columns.Item(6).Insert()

# Header row: copy the style of the existing header cells (E1) to the new F1 header.
$ws.Range("F1").Value = "Items Count"
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)

# Data rows: Items Count values (numeric 0)
$ws.Range("F2").Value = 0
$ws.Range("F3").Value = 0
$ws.Range("F4").Value = 0
